$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 315, shifting existing rows 315:366 down to 316:367
$ws.Rows("315:315").Insert()

$newRow = 315

$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"
$ws.Cells.Item($newRow, 4).Value = 44951
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100112009
$ws.Cells.Item($newRow, 7).Value = "Acelga"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 200
$ws.Cells.Item($newRow, 11).Value = 3000
$ws.Cells.Item($newRow, 12).Value = 3000
$ws.Cells.Item($newRow, 13).Value = 3000
$ws.Cells.Item($newRow, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 750
$ws.Cells.Item($newRow, 17).Value = 4
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
